$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1299.5667
$ws.Range("I98").Value = 1330.5186
$ws.Range("J98").Value = 1021
$ws.Range("K98").Value = 1330.5186
$ws.Range("L98").Value = 1021
$ws.Range("M98").Value = 167.4813999999999
$ws.Range("N98").Value = -4017

$ws.Range("H99").Value = 281.58334
$ws.Range("I99").Value = 281.8
$ws.Range("J99").Value = 280.5
$ws.Range("K99").Value = 845.4000000000001
$ws.Range("L99").Value = 841.5
$ws.Range("M99").Value = 652.5999999999999
$ws.Range("N99").Value = -3837.5

$ws.Range("H122").Value = 1299.5667
$ws.Range("I122").Value = 1330.5186
$ws.Range("J122").Value = 1021
$ws.Range("K122").Value = 3991.5558
$ws.Range("L122").Value = 3063
$ws.Range("M122").Value = -1541.5558
$ws.Range("N122").Value = -7963

$ws.Range("H132").Value = 1712.2051
$ws.Range("I132").Value = 1536.8472
$ws.Range("J132").Value = 3816.5
$ws.Range("K132").Value = 4610.5416
$ws.Range("L132").Value = 11449.5
$ws.Range("M132").Value = -2080.5416
$ws.Range("N132").Value = -16509.5

$ws.Range("H137").Value = 2333.6667
$ws.Range("I137").Value = 2172.7827
$ws.Range("J137").Value = 2862.2856
$ws.Range("K137").Value = 6518.348100000001
$ws.Range("L137").Value = 8586.856800000001
$ws.Range("M137").Value = -3968.348100000001
$ws.Range("N137").Value = -13686.8568

$ws.Range("H138").Value = 2090.2856
$ws.Range("I138").Value = 1027.5278
$ws.Range("J138").Value = 4003.25
$ws.Range("K138").Value = 3082.5834
$ws.Range("L138").Value = 12009.75
$ws.Range("M138").Value = 2057.4166
$ws.Range("N138").Value = -22289.75

$ws.Range("H141").Value = 1125.04
$ws.Range("I141").Value = 1088.5834
$ws.Range("K141").Value = 3265.7502
$ws.Range("M141").Value = 1914.2498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5720.896
$ws.Range("I32").Value = 3805.092
$ws.Range("J32").Value = 24240.334
$ws.Range("K32").Value = 3805.092
$ws.Range("L32").Value = 24240.334
$ws.Range("M32").Value = -3518.092
$ws.Range("N32").Value = -24814.334

$ws.Range("H61").Value = 22271.203
$ws.Range("I61").Value = 2943.7058
$ws.Range("J61").Value = 66080.2
$ws.Range("K61").Value = 2943.7058
$ws.Range("L61").Value = 66080.2
$ws.Range("M61").Value = -2731.7058
$ws.Range("N61").Value = -66504.2

$ws.Range("H74").Value = 140562
$ws.Range("I74").Value = 101688.11
$ws.Range("J74").Value = 228028.25
$ws.Range("K74").Value = 101688.11
$ws.Range("L74").Value = 228028.25
$ws.Range("M74").Value = -100814.11
$ws.Range("N74").Value = -229776.25

$ws.Range("H77").Value = 140562
$ws.Range("I77").Value = 101688.11
$ws.Range("J77").Value = 228028.25
$ws.Range("K77").Value = 508440.55
$ws.Range("L77").Value = 1140141.25
$ws.Range("M77").Value = -504072.55
$ws.Range("N77").Value = -1148877.25

$ws.Range("H97").Value = 830.6177
$ws.Range("I97").Value = 494.68967
$ws.Range("K97").Value = 494.68967
$ws.Range("M97").Value = 1.310330000000022

$ws.Range("H132").Value = 2581.5
$ws.Range("I132").Value = 2394.2593
$ws.Range("J132").Value = 4266.6665
$ws.Range("K132").Value = 7182.777900000001
$ws.Range("L132").Value = 12799.9995
$ws.Range("M132").Value = -4652.777900000001
$ws.Range("N132").Value = -17859.9995

$ws.Range("H136").Value = 22271.203
$ws.Range("I136").Value = 2943.7058
$ws.Range("J136").Value = 66080.2
$ws.Range("K136").Value = 8831.117400000001
$ws.Range("L136").Value = 198240.6
$ws.Range("M136").Value = -6281.117400000001
$ws.Range("N136").Value = -203340.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2433.9614
$ws.Range("I20").Value = 2520.8096
$ws.Range("K20").Value = 2520.8096
$ws.Range("M20").Value = -2273.8096

$ws.Range("H105").Value = 855.931
$ws.Range("I105").Value = 686.7083
$ws.Range("J105").Value = 1668.2
$ws.Range("K105").Value = 686.7083
$ws.Range("L105").Value = 1668.2
$ws.Range("M105").Value = 1060.2917
$ws.Range("N105").Value = -5162.2

$ws.Range("H134").Value = 2876.3257
$ws.Range("I134").Value = 1817.7878
$ws.Range("K134").Value = 5453.3634
$ws.Range("M134").Value = -2918.3634

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3103.2188
$ws.Range("I58").Value = 2806.5173
$ws.Range("J58").Value = 5971.3335
$ws.Range("K58").Value = 2806.5173
$ws.Range("L58").Value = 5971.3335
$ws.Range("M58").Value = -2603.5173
$ws.Range("N58").Value = -6377.3335

$ws.Range("H132").Value = 6384.773
$ws.Range("I132").Value = 1749
$ws.Range("J132").Value = 18746.834
$ws.Range("K132").Value = 5247
$ws.Range("L132").Value = 56240.50199999999
$ws.Range("M132").Value = -2717
$ws.Range("N132").Value = -61300.50199999999

$ws.Range("H134").Value = 2606.0876
$ws.Range("I134").Value = 2436.0962
$ws.Range("K134").Value = 7308.2886
$ws.Range("M134").Value = -4773.2886

$ws.Range("H136").Value = 3103.2188
$ws.Range("I136").Value = 2806.5173
$ws.Range("J136").Value = 5971.3335
$ws.Range("K136").Value = 8419.5519
$ws.Range("L136").Value = 17914.0005
$ws.Range("M136").Value = -5869.5519
$ws.Range("N136").Value = -23014.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3507.647
$ws.Range("I131").Value = 2209.375
$ws.Range("J131").Value = 4661.6665
$ws.Range("K131").Value = 6628.125
$ws.Range("L131").Value = 13984.9995
$ws.Range("M131").Value = -1588.125
$ws.Range("N131").Value = -24064.9995

$ws.Range("H137").Value = 3458.9333
$ws.Range("I137").Value = 1637.8
$ws.Range("J137").Value = 4369.5
$ws.Range("K137").Value = 4913.4
$ws.Range("L137").Value = 13108.5
$ws.Range("M137").Value = 186.6000000000004
$ws.Range("N137").Value = -23308.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6372.1665
$ws.Range("I70").Value = 5802.8
$ws.Range("J70").Value = 6778.857
$ws.Range("K70").Value = 5802.8
$ws.Range("L70").Value = 6778.857
$ws.Range("M70").Value = -5532.8
$ws.Range("N70").Value = -7318.857

$ws.Range("H73").Value = 6372.1665
$ws.Range("I73").Value = 5802.8
$ws.Range("J73").Value = 6778.857
$ws.Range("K73").Value = 5802.8
$ws.Range("L73").Value = 6778.857
$ws.Range("M73").Value = -4866.8
$ws.Range("N73").Value = -8650.857

$ws.Range("H102").Value = 24660.5
$ws.Range("I102").Value = 27524.05
$ws.Range("J102").Value = 15115.333
$ws.Range("K102").Value = 27524.05
$ws.Range("L102").Value = 15115.333
$ws.Range("M102").Value = -25902.05
$ws.Range("N102").Value = -18359.333

$ws.Range("H132").Value = 3924.5964
$ws.Range("I132").Value = 4135.467
$ws.Range("J132").Value = 3133.8333
$ws.Range("K132").Value = 12406.401
$ws.Range("L132").Value = 9401.499899999999
$ws.Range("M132").Value = -9876.400999999998
$ws.Range("N132").Value = -14461.4999

$ws.Range("H136").Value = 15336
$ws.Range("J136").Value = 15336
$ws.Range("L136").Value = 46008
$ws.Range("N136").Value = -51108

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1818.5
$ws.Range("I22").Value = 1737.1052
$ws.Range("J22").Value = 1921.6
$ws.Range("K22").Value = 1737.1052
$ws.Range("L22").Value = 1921.6
$ws.Range("M22").Value = -1442.1052
$ws.Range("N22").Value = -2511.6

$ws.Range("H27").Value = 1818.5
$ws.Range("I27").Value = 1737.1052
$ws.Range("J27").Value = 1921.6
$ws.Range("K27").Value = 1737.1052
$ws.Range("L27").Value = 1921.6
$ws.Range("M27").Value = -1630.1052
$ws.Range("N27").Value = -2135.6

$ws.Range("H55").Value = 1491.6207
$ws.Range("I55").Value = 947.75
$ws.Range("J55").Value = 1875.5294
$ws.Range("K55").Value = 947.75
$ws.Range("L55").Value = 1875.5294
$ws.Range("M55").Value = -774.75
$ws.Range("N55").Value = -2221.5294

$ws.Range("H132").Value = 3211.3333
$ws.Range("I132").Value = 2520.7188
$ws.Range("J132").Value = 6368.4287
$ws.Range("K132").Value = 7562.1564
$ws.Range("L132").Value = 19105.2861
$ws.Range("M132").Value = -5032.1564
$ws.Range("N132").Value = -24165.2861

$ws.Range("H136").Value = 3023.6
$ws.Range("I136").Value = 2607.0625
$ws.Range("K136").Value = 7821.1875
$ws.Range("M136").Value = -5271.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1461.0731
$ws.Range("I126").Value = 1386.1714
$ws.Range("J126").Value = 1898
$ws.Range("K126").Value = 4158.5142
$ws.Range("L126").Value = 5694
$ws.Range("M126").Value = -1688.5142
$ws.Range("N126").Value = -10634

$ws.Range("H132").Value = 4372.1875
$ws.Range("I132").Value = 3696.7693
$ws.Range("J132").Value = 7299
$ws.Range("K132").Value = 11090.3079
$ws.Range("L132").Value = 21897
$ws.Range("M132").Value = -8560.3079
$ws.Range("N132").Value = -26957
